$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Smoke Testing")

# ---- Row 1 (title) ----
$ws.Rows.Item(1).RowHeight = 96

# ---- New content for rows 3-8 (ID stays the same S-1..S-6), Title/ER/Status change ----
$ws.Range("B3").Value = "Установить приложение"
$ws.Range("C3").Value = "Приложение скачивается из Google Play, установка проходит без сбоев и ошибок."
$ws.Range("D3").Value = "Pass"

$ws.Range("B4").Value = "Запустить приложение"
$ws.Range("C4").Value = "Приложение загружается быстро или предоставляет пользователю обратную связь на экране. Материалы приложения и метаданные являются полными и точными, отражают основные возможности приложения."
$ws.Range("D4").Value = "Pass"

$ws.Range("B5").Value = "Создать список"
$ws.Range("C5").Value = "На главном экране приложения отображается созданный список."
$ws.Range("D5").Value = "Pass"

$ws.Range("B6").Value = "Добавить задачу в список"
$ws.Range("C6").Value = "В текущем списке отображается созданная задача."
$ws.Range("D6").Value = "Pass"

$ws.Range("B7").Value = "Пометить задачу, как выполненную"
$ws.Range("C7").Value = "Название выбранной задачи зачеркнуто в текущем списке, чекбокс отмечен."
$ws.Range("D7").Value = "Pass"

$ws.Range("B8").Value = "Свернуть приложение"
$ws.Range("C8").Value = "Приложение поддерживает жестовую навигацию для возврата  на главный экран."
$ws.Range("D8").Value = "Pass"

# ---- Row 9 (was S-7, now emptied) ----
$ws.Range("A9:D9").ClearContents()

# ---- Row heights ----
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 75
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 15

# ---- Status column (D3:D8) green "Pass" style ----
$statusRange = $ws.Range("D3:D8")
$statusRange.Interior.Color = 65280
$statusRange.HorizontalAlignment = -4108
$statusRange.VerticalAlignment = -4160
$statusRange.Borders.LineStyle = 1
$statusRange.Font.Name = "Corbel"

# ---- sheet view ----
$excel.ActiveWindow.Zoom = 110
$ws.Range("I3").Select()
